# AIC_Min.xlsx edit: "Creation allocation matrix for dMRWIO + Results_region"
#
# A new "Photovoltaic plants" sector row was inserted ahead of "Onshore wind
# plants" in the master Sector list (and the old trailing "Photovoltaic
# plants" entry was dropped), which shifts the 3-row Sector block -
# Onshore wind plants / Offshore wind plants / Photovoltaic plants - that is
# repeated in rows 5:7 of every yearly worksheet. The net, observable effect
# in this workbook is that each worksheet's 3-row block (label in column C,
# numeric value in column E) rotates down by one: what used to be row 7
# becomes row 5, what used to be row 5 becomes row 6, and what used to be
# row 6 becomes row 7.
#
# Apply that rotation to every worksheet (rows 5, 6 and 7; columns C and E).

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Snapshot the current (pre-rotation) label + value pairs for rows 5-7.
    $label5 = $ws.Range("C5").Value()
    $label6 = $ws.Range("C6").Value()
    $label7 = $ws.Range("C7").Value()

    $value5 = $ws.Range("E5").Value()
    $value6 = $ws.Range("E6").Value()
    $value7 = $ws.Range("E7").Value()

    # Rotate down by one: old row7 -> row5, old row5 -> row6, old row6 -> row7.
    $ws.Range("C5").Value = $label7
    $ws.Range("C6").Value = $label5
    $ws.Range("C7").Value = $label6

    $ws.Range("E5").Value = $value7
    $ws.Range("E6").Value = $value5
    $ws.Range("E7").Value = $value6
}
